$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Text = "Grupo 03`rKaggle: Grupo 03`rGithub: https://github.com/benet1one/Mineria"
for ($i=1; $i -le 3; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.LanguageID = "ca-ES"
}
